# Commit: "objects and collections export OK. abandon jdbc resultsets."
#
# This workbook is a set of `poi-tl`-style export fixtures: "Sheet2" holds a
# single-object template render, "Sheet1" holds a simple collection-of-scalars
# render ("cols"/test1/test3), and this change adds a third sheet showing a
# collection-of-objects render (a list of named colors with their hex codes),
# then leaves that new sheet selected/active.

$wb = $excel.ActiveWorkbook

# --- "Sheet2": the Dob value was re-rendered (same cell/style, new timestamp) ---
$sheetObj = $wb.Worksheets.Item("Sheet2")
$sheetObj.Range("B4").Value = 42395.67369581018

# --- "Sheet1": content ("cols" / test1 / test3) is unchanged by this commit ---
# (its tab simply stops being the selected one once Sheet3 is appended below)

# --- Append the new "Sheet3" collection-of-objects sample after "Sheet1" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet3 = $wb.Worksheets.Add($null, $lastSheet)
$sheet3.Name = "Sheet3"

$sheet3.Range("B1").Value = "Color"
$sheet3.Range("C1").Value = "Value"

$sheet3.Range("B2").Value = "red"
$sheet3.Range("C2").Value = "#f00"
$sheet3.Range("B3").Value = "green"
$sheet3.Range("C3").Value = "#0f0"
$sheet3.Range("B4").Value = "blue"
$sheet3.Range("C4").Value = "#00f"
$sheet3.Range("B5").Value = "cyan"
$sheet3.Range("C5").Value = "#0ff"
$sheet3.Range("B6").Value = "magenta"
$sheet3.Range("C6").Value = "#f0f"
$sheet3.Range("B7").Value = "yellow"
$sheet3.Range("C7").Value = "#ff0"
$sheet3.Range("B8").Value = "black"
$sheet3.Range("C8").Value = "#000"

$sheet3.Columns.Item(2).ColumnWidth = 18.7265625
$sheet3.Columns.Item(3).ColumnWidth = 19.0

# New sheet becomes the active tab, with D5 as its selected cell.
[void]$sheet3.Range("D5").Select()
